$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates: column, row, new value
# D-column (price) values may look numeric, so we force text formatting
# to preserve exact string representation (e.g. trailing zeros, multiple dots).
$updates = @(
    @('D', 2, '27.168.16'),
    @('E', 2, '  +0.60%  '),
    @('D', 3, '1.569.48'),
    @('E', 3, '  +0.58%  '),
    @('E', 4, '  +0.64%  '),
    @('D', 5, '211.60'),
    @('E', 5, '  +2.07%  '),
    @('E', 6, '  +0.49%  '),
    @('E', 7, '  +0.62%  '),
    @('D', 8, '22.01'),
    @('E', 8, '  -0.71%  '),
    @('E', 9, '  +0.26%  '),
    @('E', 10, '  +0.57%  '),
    @('D', 11, '0.0867'),
    @('E', 11, '  +0.96%  '),
    @('D', 12, '1.792.52'),
    @('E', 12, '  +0.59%  '),
    @('D', 13, '1.566.86'),
    @('E', 13, '  +0.80%  '),
    @('D', 14, '3.79'),
    @('E', 14, '  +0.60%  '),
    @('D', 15, '0.519'),
    @('E', 15, '  -0.44%  '),
    @('B', 16, 'Litecoin'),
    @('C', 16, 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'),
    @('D', 16, '62.34'),
    @('E', 16, '  +0.49%  '),
    @('B', 17, 'WrappedBTC'),
    @('C', 17, 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'),
    @('D', 17, '27.155.16'),
    @('E', 17, '  +0.59%  '),
    @('D', 18, '0.0₃0703'),
    @('E', 18, '  -0.68%  '),
    @('D', 19, '215.94'),
    @('E', 19, '  -0.66%  '),
    @('E', 20, '  +0.67%  '),
    @('E', 21, '  +0.53%  '),
    @('E', 22, '  +1.18%  '),
    @('E', 23, '  -0.20%  '),
    @('E', 24, '  +0.66%  '),
    @('D', 25, '153.95'),
    @('E', 25, '  +0.52%  '),
    @('D', 26, '6.65'),
    @('E', 26, '  +0.04%  '),
    @('D', 27, '15.08'),
    @('E', 27, '  +0.30%  '),
    @('D', 28, '0.106'),
    @('E', 28, '  +1.68%  '),
    @('E', 29, '  +0.37%  '),
    @('E', 30, '  +2.35%  '),
    @('D', 31, '0.0473'),
    @('E', 31, '  +0.93%  '),
    @('E', 32, '  +0.11%  '),
    @('D', 33, '3.18'),
    @('E', 33, '  +2.40%  '),
    @('D', 34, '1.449.92'),
    @('E', 34, '  +1.93%  '),
    @('E', 35, '  +5.86%  '),
    @('E', 36, '  +0.33%  '),
    @('E', 37, '  +1.66%  '),
    @('D', 38, '0.0167'),
    @('E', 38, '  +1.07%  '),
    @('E', 39, '  +0.55%  '),
    @('E', 40, '  +2.73%  '),
    @('E', 41, '  -0.03%  '),
    @('E', 42, '  +0.57%  '),
    @('D', 43, '2.34'),
    @('E', 43, '  +0.81%  '),
    @('E', 44, '  -1.28%  '),
    @('D', 45, '64.67'),
    @('E', 45, '  -0.34%  '),
    @('E', 46, '  -1.14%  '),
    @('D', 47, '1.705.22'),
    @('E', 47, '  +0.59%  '),
    @('D', 48, '86.00'),
    @('E', 48, '  -1.67%  '),
    @('D', 49, '0.0₆0103'),
    @('E', 49, '  +2.16%  '),
    @('D', 50, '0.0520'),
    @('E', 50, '  -0.03%  '),
    @('D', 51, '0.0959'),
    @('E', 51, '  -0.04%  ')
)

foreach ($u in $updates) {
    $col = $u[0]
    $row = $u[1]
    $val = $u[2]
    $cell = $ws.Range("$col$row")
    if ($col -eq "D") {
        # Force text number format so numeric-looking strings (e.g. "211.60")
        # are not auto-converted into Excel numbers, losing formatting.
        $cell.NumberFormat = "@"
        $cell.Value = $val
        $cell.Style = "Normal"
    } else {
        $cell.Value = $val
    }
}
